# Contest 43 RR vs RCB.
# Fill in the raw score entries for row 55 (Contest 43, "RR vs RCB").
# The D/G/J/M/P/S/V/Y/AB columns already hold VLOOKUP/RANK formulas that
# derive their result from the adjacent raw-score cell, so only the raw
# score cells (E/H/K/N/Q/T/W/Z/AC) need to be populated; Excel recalculates
# the formula cells (and the row 60 SUM totals) automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E55").Value  = 20
$ws.Range("H55").Value  = 40
$ws.Range("K55").Value  = 60
$ws.Range("N55").Value  = 100
$ws.Range("Q55").Value  = 70
$ws.Range("T55").Value  = 50
$ws.Range("W55").Value  = 0
$ws.Range("Z55").Value  = 80
$ws.Range("AC55").Value = 30
